$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "001" to "004" (must stay text, not become numeric 4)
$ws.Range("J2").Value = "'004"

# Report date update
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Updated numeric figures
$ws.Range("O2").Value = 33823349.92
$ws.Range("P2").Value = 273.857785228
$ws.Range("Q2").Value = 208212588.59
$ws.Range("R2").Value = 1685.8365153868
$ws.Range("S2").Value = 87655458.06
$ws.Range("T2").Value = 709.7206416346
$ws.Range("U2").Value = 4561392.31
$ws.Range("V2").Value = 36.9322612493
$ws.Range("W2").Value = 588716.66
$ws.Range("X2").Value = 4.7666668445
$ws.Range("Y2").Value = 8027424.35
$ws.Range("Z2").Value = 64.9957103237
$ws.Range("AA2").Value = -25364228.93
$ws.Range("AB2").Value = -205.3667532996
$ws.Range("AC2").Value = 12350698.7
$ws.Range("AD2").Value = 165.3346002479
